{"js": "// The document contains a single 20x5 table of simple arithmetic\n// problems (e.g. \"62-27=\"). The commit replaces the text of every\n// cell with a new arithmetic problem, cell-for-cell, in row-major\n// (reading) order. The surrounding formatting/structure is untouched.\n\nconst newValues = [\n  [\"98-29=\", \"70-23=\", \"10-1=\", \"3+78=\", \"43-29=\"],\n  [\"28+57=\", \"57-18=\", \"29+39=\", \"49+3=\", \"38+38=\"],\n  [\"84-19=\", \"82-65=\", \"7+47=\", \"12+9=\", \"62-5=\"],\n  [\"9+6=\", \"83-64=\", \"34+48=\", \"26+16=\", \"59+15=\"],\n  [\"38+48=\", \"79+18=\", \"30-19=\", \"23+68=\", \"16+65=\"],\n  [\"40-26=\", \"9+38=\", \"77+15=\", \"89+9=\", \"8+65=\"],\n  [\"8+26=\", \"78+3=\", \"60-33=\", \"76+16=\", \"81-5=\"],\n  [\"26+67=\", \"58+17=\", \"90-26=\", \"91-44=\", \"41-29=\"],\n  [\"9+37=\", \"71-17=\", \"7+64=\", \"71-24=\", \"10-6=\"],\n  [\"53-26=\", \"27+45=\", \"62-16=\", \"21-7=\", \"7+44=\"],\n  [\"3+78=\", \"40-22=\", \"39+44=\", \"53-44=\", \"56+17=\"],\n  [\"57-29=\", \"68+3=\", \"32-4=\", \"34+57=\", \"56+28=\"],\n  [\"68+19=\", \"85-26=\", \"5+77=\", \"27+6=\", \"37+36=\"],\n  [\"49+8=\", \"17+56=\", \"56+17=\", \"14+67=\", \"85-36=\"],\n  [\"19+35=\", \"88+3=\", \"39+12=\", \"94-68=\", \"9+19=\"],\n  [\"68+24=\", \"31-24=\", \"44+19=\", \"94-39=\", \"48+19=\"],\n  [\"34+29=\", \"34+27=\", \"91-32=\", \"94-79=\", \"9+52=\"],\n  [\"9+46=\", \"83+9=\", \"14+59=\", \"95-66=\", \"47+28=\"],\n  [\"17+27=\", \"84+9=\", \"41-32=\", \"13+69=\", \"93-46=\"],\n  [\"72-7=\", \"92-47=\", \"84-77=\", \"92-53=\", \"75-68=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Overwrite every cell's text in place (keeps existing run/paragraph\n// formatting since Word only swaps the text runs, not the cells).\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single 20x5 table of simple arithmetic\n# problems (e.g. \"62-27=\"). The commit replaces the text of every\n# cell with a new arithmetic problem, cell-for-cell, in row-major\n# (reading) order. The surrounding formatting/structure is untouched.\n\n$newValues = @(\n  @(\"98-29=\",\"70-23=\",\"10-1=\",\"3+78=\",\"43-29=\"),\n  @(\"28+57=\",\"57-18=\",\"29+39=\",\"49+3=\",\"38+38=\"),\n  @(\"84-19=\",\"82-65=\",\"7+47=\",\"12+9=\",\"62-5=\"),\n  @(\"9+6=\",\"83-64=\",\"34+48=\",\"26+16=\",\"59+15=\"),\n  @(\"38+48=\",\"79+18=\",\"30-19=\",\"23+68=\",\"16+65=\"),\n  @(\"40-26=\",\"9+38=\",\"77+15=\",\"89+9=\",\"8+65=\"),\n  @(\"8+26=\",\"78+3=\",\"60-33=\",\"76+16=\",\"81-5=\"),\n  @(\"26+67=\",\"58+17=\",\"90-26=\",\"91-44=\",\"41-29=\"),\n  @(\"9+37=\",\"71-17=\",\"7+64=\",\"71-24=\",\"10-6=\"),\n  @(\"53-26=\",\"27+45=\",\"62-16=\",\"21-7=\",\"7+44=\"),\n  @(\"3+78=\",\"40-22=\",\"39+44=\",\"53-44=\",\"56+17=\"),\n  @(\"57-29=\",\"68+3=\",\"32-4=\",\"34+57=\",\"56+28=\"),\n  @(\"68+19=\",\"85-26=\",\"5+77=\",\"27+6=\",\"37+36=\"),\n  @(\"49+8=\",\"17+56=\",\"56+17=\",\"14+67=\",\"85-36=\"),\n  @(\"19+35=\",\"88+3=\",\"39+12=\",\"94-68=\",\"9+19=\"),\n  @(\"68+24=\",\"31-24=\",\"44+19=\",\"94-39=\",\"48+19=\"),\n  @(\"34+29=\",\"34+27=\",\"91-32=\",\"94-79=\",\"9+52=\"),\n  @(\"9+46=\",\"83+9=\",\"14+59=\",\"95-66=\",\"47+28=\"),\n  @(\"17+27=\",\"84+9=\",\"41-32=\",\"13+69=\",\"93-46=\"),\n  @(\"72-7=\",\"92-47=\",\"84-77=\",\"92-53=\",\"75-68=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n\n"}
